# Update NATMI LR-pairs output values (Anxa1-Egfr) following Dr Hou's advice.
# The ligand/receptor-expressing cell counts (columns E and K) change from 1 to 3,
# and the dependent expression / specificity statistics (G,H,I,J,M,N,O,P,Q,R,S,T)
# are recomputed accordingly for data rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "E" = 3; "G" = 59.68504933333333; "H" = 179.055148; "I" = 0.207862575863973; "J" = 0.2078625758639731; "K" = 3; "M" = 2.08532; "N" = 6.25596; "O" = 0.01753772176136817; "P" = 0.01753772176136816; "Q" = 124.4624270757866; "R" = 1120.16184368208; "S" = 0.003645436020103641; "T" = 0.003645436020103641 }
    3 = @{ "E" = 3; "G" = 59.68504933333333; "H" = 179.055148; "I" = 0.207862575863973; "J" = 0.2078625758639731; "K" = 3; "M" = 101.898173; "N" = 305.694519; "O" = 0.8569724579756384; "P" = 0.8569724579756383; "Q" = 6081.797482481535; "R" = 54736.17734233382; "S" = 0.1781325025592966; "T" = 0.1781325025592966 }
    4 = @{ "E" = 3; "G" = 59.68504933333333; "H" = 179.055148; "I" = 0.207862575863973; "J" = 0.2078625758639731; "K" = 3; "M" = 14.921347; "N" = 44.76404100000001; "O" = 0.1254898202629935; "P" = 0.1254898202629935; "Q" = 890.5813318147855; "R" = 8015.231986333069; "S" = 0.02608463728457283; "T" = 0.02608463728457283 }
    5 = @{ "E" = 3; "G" = 186.4134573333333; "H" = 559.240372; "I" = 0.6492141976897894; "J" = 0.6492141976897894; "K" = 3; "M" = 2.08532; "N" = 6.25596; "O" = 0.01753772176136817; "P" = 0.01753772176136816; "Q" = 388.7317108463466; "R" = 3498.58539761712; "S" = 0.01138573796261339; "T" = 0.01138573796261339 }
    6 = @{ "E" = 3; "G" = 186.4134573333333; "H" = 559.240372; "I" = 0.6492141976897894; "J" = 0.6492141976897894; "K" = 3; "M" = 101.898173; "N" = 305.694519; "O" = 0.8569724579756384; "P" = 0.8569724579756383; "Q" = 18995.19072488012; "R" = 170956.7165239211; "S" = 0.5563586867469008; "T" = 0.5563586867469007 }
    7 = @{ "E" = 3; "G" = 186.4134573333333; "H" = 559.240372; "I" = 0.6492141976897894; "J" = 0.6492141976897894; "K" = 3; "M" = 14.921347; "N" = 44.76404100000001; "O" = 0.1254898202629935; "P" = 0.1254898202629935; "Q" = 2781.539882340362; "R" = 25033.85894106326; "S" = 0.08146977298027522; "T" = 0.08146977298027519 }
    8 = @{ "E" = 3; "G" = 41.03855533333333; "H" = 123.115666; "I" = 0.1429232264462375; "J" = 0.1429232264462375; "K" = 3; "M" = 2.08532; "N" = 6.25596; "O" = 0.01753772176136817; "P" = 0.01753772176136816; "Q" = 85.57852020770666; "R" = 770.2066818693601; "S" = 0.002506547778651129; "T" = 0.002506547778651129 }
    9 = @{ "E" = 3; "G" = 41.03855533333333; "H" = 123.115666; "I" = 0.1429232264462375; "J" = 0.1429232264462375; "K" = 3; "M" = 101.898173; "N" = 305.694519; "O" = 0.8569724579756384; "P" = 0.8569724579756383; "Q" = 4181.753811026073; "R" = 37635.78429923466; "S" = 0.1224812686694409; "T" = 0.1224812686694409 }
    10 = @{ "E" = 3; "G" = 41.03855533333333; "H" = 123.115666; "I" = 0.1429232264462375; "J" = 0.1429232264462375; "K" = 3; "M" = 14.921347; "N" = 44.76404100000001; "O" = 0.1254898202629935; "P" = 0.1254898202629935; "Q" = 612.3505245073675; "R" = 5511.154720566307; "S" = 0.01793540999814546; "T" = 0.01793540999814546 }
}

foreach ($rowKey in $updates.Keys) {
    $rowData = $updates[$rowKey]
    foreach ($colKey in $rowData.Keys) {
        $ws.Range("$colKey$rowKey").Value = $rowData[$colKey]
    }
}
